# Apply the commit's change:
#  - Remove the "Identifier" property row from the Metadata sheet
#    (row 3: "Identifier" / "http://interopsante.org/fhir/ValueSet#urn:oid:2.16.840.1.113883.2.8.1.3.4")
#  - Update the "Date" property value to the new timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Delete the whole "Identifier" row (row 3) - shifts everything below up by one.
$ws.Rows("3:3").Delete()

# Find the "Date" row (was row 9, now row 8 after the deletion above) and update its value.
$dateCell = $ws.Range("A8")
if ($dateCell.Value2 -eq "Date") {
    $ws.Range("B8").Value = "2025-12-02T20:15:58+00:00"
} else {
    # Fallback: search for the row containing "Date" in column A
    for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
        if ($ws.Cells.Item($r, 1).Value2 -eq "Date") {
            $ws.Cells.Item($r, 2).Value = "2025-12-02T20:15:58+00:00"
            break
        }
    }
}
